$d = $word.ActiveDocument

# 1. Merge "Idade: 2" + "6" into a single run "Idade: 26"
$p = $d.Paragraphs(3)
$rng = $d.Range($p.Range.Start, $p.Range.End - 1)
$rng.InsertXML('<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>Idade: 26</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>')

# 2. Add "Java, Javascript, C#, " run before the Nodejs... run
$p = $d.Paragraphs(18)
$rng = $d.Range($p.Range.Start, $p.Range.End - 1)
$rng.InsertXML('<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:rPr/><w:t xml:space="preserve">Java, Javascript, C#, </w:t></w:r><w:r><w:rPr/><w:t>Nodejs, es6/es7, PHP, HTML5. CSS, MongoDB, My-SQL, PL-SQL</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>')

# 3. Move the bookmark from the "Frameworks" paragraph to the "React..." paragraph
$p19 = $d.Paragraphs(19)
$p22 = $d.Paragraphs(22)
$rng = $d.Range($p19.Range.Start, $p22.Range.End)
$rng.InsertXML('<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="TextBody"/><w:rPr><w:b/><w:b/><w:bCs/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>Frameworks</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="TextBody"/><w:rPr><w:b w:val="false"/><w:b w:val="false"/><w:bCs w:val="false"/></w:rPr></w:pPr><w:bookmarkStart w:id="0" w:name="__DdeLink__31_290923235"/><w:r><w:rPr><w:b w:val="false"/><w:bCs w:val="false"/></w:rPr><w:t>React, Redux, Axios, Express, Mongoose</w:t></w:r><w:bookmarkEnd w:id="0"/><w:r><w:rPr><w:b w:val="false"/><w:bCs w:val="false"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="TextBody"/><w:rPr><w:b/><w:b/><w:bCs/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>Game Development</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="TextBody"/><w:rPr><w:b w:val="false"/><w:b w:val="false"/><w:bCs w:val="false"/></w:rPr></w:pPr><w:r><w:rPr><w:b w:val="false"/><w:bCs w:val="false"/></w:rPr><w:t>Unity, C#</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>')

# 4. Merge "Github, " + "Amazon Web Services, Heroku, Google Play Services" into one run
$p = $d.Paragraphs(24)
$rng = $d.Range($p.Range.Start, $p.Range.End - 1)
$rng.InsertXML('<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:rPr/><w:t>Github, Amazon Web Services, Heroku, Google Play Services</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>')

# 5. Merge "Applications " + "& tools" into one run
$p = $d.Paragraphs(25)
$rng = $d.Range($p.Range.Start, $p.Range.End - 1)
$rng.InsertXML('<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:rPr><w:b/></w:rPr><w:t>Applications &amp; tools</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>')

# 6. Merge "Adobe Photoshop, ... Android Studio, " + "Discord" into one run
$p = $d.Paragraphs(26)
$rng = $d.Range($p.Range.Start, $p.Range.End - 1)
$rng.InsertXML('<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:rPr/><w:t>Adobe Photoshop,  Gimp, Atom, Android Studio, Discord</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>')

